# Update market price / profit figures pulled by the scheduled runner.
# Applies literal value updates to columns H-N for specific Leve rows
# across several job sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4966.8096
$ws.Range("I64").Value = 3956.6667
$ws.Range("J64").Value = 6313.6665
$ws.Range("K64").Value = 3956.6667
$ws.Range("L64").Value = 6313.6665
$ws.Range("M64").Value = -3708.6667
$ws.Range("N64").Value = -6809.6665

$ws.Range("H67").Value = 4966.8096
$ws.Range("I67").Value = 3956.6667
$ws.Range("J67").Value = 6313.6665
$ws.Range("K67").Value = 3956.6667
$ws.Range("L67").Value = 6313.6665
$ws.Range("M67").Value = -3098.6667
$ws.Range("N67").Value = -8029.6665

$ws.Range("H76").Value = 6947631.5
$ws.Range("I76").Value = 11113811
$ws.Range("K76").Value = 11113811
$ws.Range("M76").Value = -11113496

$ws.Range("H79").Value = 6947631.5
$ws.Range("I79").Value = 11113811
$ws.Range("K79").Value = 11113811
$ws.Range("M79").Value = -11112719

$ws.Range("H133").Value = 45936
$ws.Range("J133").Value = 45936
$ws.Range("L133").Value = 45936
$ws.Range("N133").Value = -56056

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3157.4167
$ws.Range("I61").Value = 2989.9092
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2989.9092
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2777.9092
$ws.Range("N61").Value = -5424

$ws.Range("H133").Value = 45587
$ws.Range("J133").Value = 45587
$ws.Range("L133").Value = 45587
$ws.Range("N133").Value = -50647

$ws.Range("H136").Value = 3157.4167
$ws.Range("I136").Value = 2989.9092
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 8969.7276
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -6419.7276
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3198.8262
$ws.Range("I105").Value = 2988.9395
$ws.Range("J105").Value = 3731.6155
$ws.Range("K105").Value = 2988.9395
$ws.Range("L105").Value = 3731.6155
$ws.Range("M105").Value = -1241.9395
$ws.Range("N105").Value = -7225.6155

$ws.Range("H133").Value = 43000
$ws.Range("J133").Value = 43000
$ws.Range("L133").Value = 43000
$ws.Range("N133").Value = -53120

$ws.Range("H139").Value = 94560
$ws.Range("J139").Value = 94560
$ws.Range("L139").Value = 94560
$ws.Range("N139").Value = -104840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4254.1025
$ws.Range("I31").Value = 1970.2106
$ws.Range("J31").Value = 6423.8
$ws.Range("K31").Value = 1970.2106
$ws.Range("L31").Value = 6423.8
$ws.Range("M31").Value = -1675.2106
$ws.Range("N31").Value = -7013.8

$ws.Range("H34").Value = 4254.1025
$ws.Range("I34").Value = 1970.2106
$ws.Range("J34").Value = 6423.8
$ws.Range("K34").Value = 1970.2106
$ws.Range("L34").Value = 6423.8
$ws.Range("M34").Value = -1768.2106
$ws.Range("N34").Value = -6827.8

$ws.Range("H58").Value = 1561.8
$ws.Range("I58").Value = 1115.8462
$ws.Range("J58").Value = 4460.5
$ws.Range("K58").Value = 1115.8462
$ws.Range("L58").Value = 4460.5
$ws.Range("M58").Value = -912.8462
$ws.Range("N58").Value = -4866.5

$ws.Range("H62").Value = 17643.47
$ws.Range("I62").Value = 20933.77
$ws.Range("J62").Value = 6950
$ws.Range("K62").Value = 20933.77
$ws.Range("L62").Value = 6950
$ws.Range("M62").Value = -20309.77
$ws.Range("N62").Value = -8198

$ws.Range("H65").Value = 17643.47
$ws.Range("I65").Value = 20933.77
$ws.Range("J65").Value = 6950
$ws.Range("K65").Value = 104668.85
$ws.Range("L65").Value = 34750
$ws.Range("M65").Value = -101548.85
$ws.Range("N65").Value = -40990

$ws.Range("H136").Value = 1561.8
$ws.Range("I136").Value = 1115.8462
$ws.Range("J136").Value = 4460.5
$ws.Range("K136").Value = 3347.5386
$ws.Range("L136").Value = 13381.5
$ws.Range("M136").Value = -797.5385999999999
$ws.Range("N136").Value = -18481.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 516.6799999999999
$ws.Range("I107").Value = 525.2143
$ws.Range("J107").Value = 505.81818
$ws.Range("K107").Value = 1575.6429
$ws.Range("L107").Value = 1517.45454
$ws.Range("M107").Value = 344.3571000000002
$ws.Range("N107").Value = -5357.45454

$ws.Range("H113").Value = 684.46155
$ws.Range("I113").Value = 590
$ws.Range("K113").Value = 1770
$ws.Range("M113").Value = 400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6102.1665
$ws.Range("I70").Value = 6375.6665
$ws.Range("J70").Value = 5281.6665
$ws.Range("K70").Value = 6375.6665
$ws.Range("L70").Value = 5281.6665
$ws.Range("M70").Value = -6105.6665
$ws.Range("N70").Value = -5821.6665

$ws.Range("H73").Value = 6102.1665
$ws.Range("I73").Value = 6375.6665
$ws.Range("J73").Value = 5281.6665
$ws.Range("K73").Value = 6375.6665
$ws.Range("L73").Value = 5281.6665
$ws.Range("M73").Value = -5439.6665
$ws.Range("N73").Value = -7153.6665

$ws.Range("H80").Value = 2963.4211
$ws.Range("I80").Value = 2700.3572
$ws.Range("J80").Value = 3700
$ws.Range("K80").Value = 2700.3572
$ws.Range("L80").Value = 3700
$ws.Range("M80").Value = -1702.3572
$ws.Range("N80").Value = -5696

$ws.Range("H83").Value = 2963.4211
$ws.Range("I83").Value = 2700.3572
$ws.Range("J83").Value = 3700
$ws.Range("K83").Value = 13501.786
$ws.Range("L83").Value = 18500
$ws.Range("M83").Value = -8509.786
$ws.Range("N83").Value = -28484

$ws.Range("H138").Value = 78000
$ws.Range("J138").Value = 78000
$ws.Range("L138").Value = 78000
$ws.Range("N138").Value = -88280

$ws.Range("H139").Value = 47666.668
$ws.Range("J139").Value = 47666.668
$ws.Range("L139").Value = 47666.668
$ws.Range("N139").Value = -57946.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1495.8462
$ws.Range("I96").Value = 1278.2222
$ws.Range("J96").Value = 1985.5
$ws.Range("K96").Value = 1278.2222
$ws.Range("L96").Value = 1985.5
$ws.Range("M96").Value = 94.77780000000007
$ws.Range("N96").Value = -4731.5

$ws.Range("H100").Value = 1500
$ws.Range("I100").Value = 1500
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459
